$d = $word.ActiveDocument

# 1. "in an academic setting" -> "in an academic research setting"
$d.Content.Find.Execute("an academic setting", $true, $false, $false, $false, $false,
                         $true, 1, $false, "an academic research setting", 2)

# 2. "seeking an internship with ITW for the summer of 2018" -> "seeking an internship for the summer of 2018"
$d.Content.Find.Execute("seeking an internship with ITW for the summer of 2018", $true, $false, $false, $false, $false,
                         $true, 1, $false, "seeking an internship for the summer of 2018", 2)
